# Update "想去人数" (want-to-go count) figures in the F column, and the
# sold-out status text in G25, on both the "展览" and "全部类型" sheets.
# (These two sheets list the same events and are kept in sync.)

$wb = $excel.ActiveWorkbook

# Column F (想去人数) new values, keyed by row number. Both sheets converge
# on the same final numbers even though a couple of rows (32, 34) started
# from slightly different counts on each sheet.
$fUpdates = @{
    2  = 143
    4  = 65
    5  = 99
    6  = 136
    7  = 1315
    8  = 1558
    9  = 343
    10 = 423
    12 = 172
    14 = 71
    15 = 113
    16 = 277
    17 = 317
    18 = 330
    19 = 1759
    20 = 72
    22 = 179
    24 = 304
    25 = 345
    26 = 4230
    28 = 286
    29 = 1111
    32 = 606
    34 = 301
    35 = 49
    36 = 151
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    foreach ($row in $fUpdates.Keys) {
        $ws.Cells.Item($row, 6).Value = $fUpdates[$row]
    }

    # G25: ticket status text changed from "暂时售罄" (temporarily sold out)
    # to "已售罄" (sold out).
    $ws.Cells.Item(25, 7).Value = "已售罄"
}
